$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X2").NumberFormat = "@"
$ws.Range("X2").Value = "87"
$ws.Range("X2").Style = "Normal"
$ws.Range("Y2").Value = "3/27;4;4/301;7;8;Alz 1;12;Ds 114;18;M 347;210;19;20;24;26;26/002;26/010;28;29;29 T;32;33;34/50;35;45;77/006;83;84;90;91;93;104;107;109;291/407;352;1019;1050;355;155;361;37;385;Tr 5;468;70;517;700;705;715;739;Zi 4;Hil 2;Bo 133;Ch 1;Brt 1;Ilm 1;115;542;841;Kt 1;897;959;Wis 1;J 59;Bs 68;519;Ki 95;Gat 1;136;1156;Dm 19;944;Dü 62;836;L 322;51;Freis 2;706;L 189;1373;960/1;960/3;Kob 7;1251;287;L 97"
